$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.202.03"
$ws.Range("E2").Value = "  -1.60%  "
$ws.Range("D3").Value = "1.854.97"
$ws.Range("E3").Value = "  -2.32%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'233.88"
$ws.Range("E5").Value = "  -2.22%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "'0.4702"
$ws.Range("E7").Value = "  -1.99%  "
$ws.Range("D8").Value = "'0.2722"
$ws.Range("E8").Value = "  -3.93%  "
$ws.Range("D9").Value = "'0.06391"
$ws.Range("E9").Value = "  -2.28%  "
$ws.Range("D10").Value = "1.842.31"
$ws.Range("E10").Value = "  -3.24%  "
$ws.Range("D11").Value = "'0.07442"
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("D12").Value = "'16.27"
$ws.Range("E12").Value = "  -2.43%  "
$ws.Range("E13").Value = "  -2.74%  "
$ws.Range("D14").Value = "'85.22"
$ws.Range("E14").Value = "  -3.12%  "
$ws.Range("D15").Value = "'0.6298"
$ws.Range("E15").Value = "  -5.67%  "
$ws.Range("D16").Value = "30.168.62"
$ws.Range("E16").Value = "  -1.63%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").Value = "'229.66"
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("D19").Value = "'12.67"
$ws.Range("E19").Value = "  -4.79%  "
$ws.Range("D20").Value = "'0.000007345"
$ws.Range("E20").Value = "  -3.33%  "
$ws.Range("D21").Value = "2.102.02"
$ws.Range("E21").Value = "  -4.38%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").Value = "'4.965"
$ws.Range("E23").Value = "  -6.48%  "
$ws.Range("D24").Value = "'5.967"
$ws.Range("E24").Value = "  -4.06%  "
$ws.Range("D25").Value = "'9.269"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").Value = "'166.14"
$ws.Range("E26").Value = "  -1.78%  "
$ws.Range("D27").Value = "'17.88"
$ws.Range("E27").Value = "  -3.74%  "
$ws.Range("D28").Value = "'1.876"
$ws.Range("E28").Value = "  -3.89%  "
$ws.Range("D29").Value = "'0.1042"
$ws.Range("E29").Value = "  +5.53%  "
$ws.Range("D30").Value = "'1.390"
$ws.Range("E30").Value = "  -0.96%  "
$ws.Range("D31").Value = "'4.134"
$ws.Range("E31").Value = "  -4.86%  "
$ws.Range("D32").Value = "'3.896"
$ws.Range("E32").Value = "  -3.03%  "
$ws.Range("D33").Value = "'0.04914"
$ws.Range("E33").Value = "  -2.80%  "
$ws.Range("D34").Value = "'1.161"
$ws.Range("E34").Value = "  -5.12%  "
$ws.Range("D35").Value = "'0.7184"
$ws.Range("E35").Value = "  -4.54%  "
$ws.Range("E36").Value = "  -0.32%  "
$ws.Range("D37").Value = "'2.705"
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("D39").Value = "'2.641"
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("D40").Value = "'0.9109"
$ws.Range("E40").Value = "  -1.03%  "
$ws.Range("D41").Value = "'1.959"
$ws.Range("E41").Value = "  -5.78%  "
$ws.Range("D42").Value = "'105.38"
$ws.Range("E42").Value = "  -1.77%  "
$ws.Range("D43").Value = "'0.9996"
$ws.Range("E43").Value = "  -0.51%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'5.551"
$ws.Range("E44").Value = "  -4.91%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "'0.4081"
$ws.Range("E45").Value = "  -4.90%  "
$ws.Range("D46").Value = "'7.081"
$ws.Range("E46").Value = "  -4.19%  "
$ws.Range("D47").Value = "'60.48"
$ws.Range("E47").Value = "  -6.27%  "
$ws.Range("D48").Value = "'0.1193"
$ws.Range("E48").Value = "  -6.20%  "
$ws.Range("D49").Value = "'8.679"
$ws.Range("E49").Value = "  -3.11%  "
$ws.Range("D50").Value = "'33.25"
$ws.Range("E50").Value = "  -1.75%  "
$ws.Range("D51").Value = "'1.395"
$ws.Range("E51").Value = "  -5.80%  "
